$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update existing correlation / p-value matrices with refreshed numbers
#    (re-run of the analysis with an extra variable folded in).
# ---------------------------------------------------------------------------

# Rho (Leptin vs Food)
$ws = $wb.Worksheets.Item("Rho")
$ws.Range("C2").Value = 0.72013378679852924
$ws.Range("B3").Value = 0.72013378679852924

# Pval (Leptin vs Food)
$ws = $wb.Worksheets.Item("Pval")
$ws.Columns.Item(3).ColumnWidth = 12.75
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.043942430227524053
$ws.Range("B3").Value = 0.043942430227524053
$ws.Range("C3").Value = 0

# Rho_G (Leptin, Food, Genotype)
$ws = $wb.Worksheets.Item("Rho_G")
$ws.Columns.Item(4).ColumnWidth = 11.75
$ws.Range("C2").Value = 0.54979415160497069
$ws.Range("D2").Value = 0.45947229012478058
$ws.Range("B3").Value = 0.54979415160497069
$ws.Range("D3").Value = 0.17420586775969896
$ws.Range("B4").Value = 0.45947229012478058
$ws.Range("C4").Value = 0.17420586775969896

# Pval_G (Leptin, Food, Genotype)
$ws = $wb.Worksheets.Item("Pval_G")
$ws.Columns.Item(4).ColumnWidth = 11.75
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.20105599132828034
$ws.Range("D2").Value = 0.29963560425619196
$ws.Range("B3").Value = 0.20105599132828034
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.70872595954746398
$ws.Range("B4").Value = 0.29963560425619196
$ws.Range("C4").Value = 0.70872595954746398
$ws.Range("D4").Value = 0

# Rho_C (Leptin, Food, Calories)
$ws = $wb.Worksheets.Item("Rho_C")
$ws.Columns.Item(4).ColumnWidth = 11.75
$ws.Range("C2").Value = 0.77779417876341361
$ws.Range("D2").Value = 0.49118618773282774
$ws.Range("B3").Value = 0.77779417876341361
$ws.Range("D3").Value = -0.77241382698286709
$ws.Range("B4").Value = 0.49118618773282774
$ws.Range("C4").Value = -0.77241382698286709

# Pval_C (Leptin, Food, Calories)
$ws = $wb.Worksheets.Item("Pval_C")
$ws.Columns.Item(4).ColumnWidth = 12.75
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.039498896219271794
$ws.Range("D2").Value = 0.26297393395023527
$ws.Range("B3").Value = 0.039498896219271794
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.041802829022008131
$ws.Range("B4").Value = 0.26297393395023527
$ws.Range("C4").Value = 0.041802829022008131
$ws.Range("D4").Value = 0

# ---------------------------------------------------------------------------
# 2) Add the two new sheets (Rho_CG, Pval_CG) covering Leptin / Food /
#    Calories / Genotype / BMI, appended right after Pval_C.
# ---------------------------------------------------------------------------

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$rhoCG = $wb.Worksheets.Add($null, $lastSheet)
$rhoCG.Name = "Rho_CG"

$pvalCG = $wb.Worksheets.Add($null, $rhoCG)
$pvalCG.Name = "Pval_CG"

# --- Rho_CG ----------------------------------------------------------------
$ws = $rhoCG
$ws.Columns.Item(1).ColumnWidth = 7.92
$ws.Columns.Item(2).ColumnWidth = 11.92
$ws.Columns.Item(3).ColumnWidth = 12.92
$ws.Columns.Item(4).ColumnWidth = 12.92
$ws.Columns.Item(5).ColumnWidth = 12.92
$ws.Columns.Item(6).ColumnWidth = 12.92

$ws.Range("A1:F1").NumberFormat = "@"
$ws.Range("A2:A6").NumberFormat = "@"

$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Leptin"
$ws.Range("C1").Value = "Food"
$ws.Range("D1").Value = "Calories"
$ws.Range("E1").Value = "Genotype"
$ws.Range("F1").Value = "BMI"

$ws.Range("A2").Value = "Leptin"
$ws.Range("A3").Value = "Food"
$ws.Range("A4").Value = "Calories"
$ws.Range("A5").Value = "Genotype"
$ws.Range("A6").Value = "BMI"

$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 0.84968511906404609
$ws.Range("D2").Value = 0.78246180092465345
$ws.Range("E2").Value = 0.75802313541041821
$ws.Range("F2").Value = 0.38073215791696413

$ws.Range("B3").Value = 0.84968511906404609
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = -0.85530743861095571
$ws.Range("E3").Value = -0.54928224222186928
$ws.Range("F3").Value = -0.41062505566534957

$ws.Range("B4").Value = 0.78246180092465345
$ws.Range("C4").Value = -0.85530743861095571
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = -0.71425135970541986
$ws.Range("F4").Value = -0.38925105153442541

$ws.Range("B5").Value = 0.75802313541041821
$ws.Range("C5").Value = -0.54928224222186928
$ws.Range("D5").Value = -0.71425135970541986
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = -0.28569356852692812

$ws.Range("B6").Value = 0.38073215791696413
$ws.Range("C6").Value = -0.41062505566534957
$ws.Range("D6").Value = -0.38925105153442541
$ws.Range("E6").Value = -0.28569356852692812
$ws.Range("F6").Value = 1

# --- Pval_CG -----------------------------------------------------------
$ws = $pvalCG
$ws.Columns.Item(1).ColumnWidth = 7.92
$ws.Columns.Item(2).ColumnWidth = 12.92
$ws.Columns.Item(3).ColumnWidth = 12.92
$ws.Columns.Item(4).ColumnWidth = 12.92
$ws.Columns.Item(5).ColumnWidth = 11.92
$ws.Columns.Item(6).ColumnWidth = 11.92

$ws.Range("A1:F1").NumberFormat = "@"
$ws.Range("A2:A6").NumberFormat = "@"

$ws.Range("A1").Value = "Row"
$ws.Range("B1").Value = "Leptin"
$ws.Range("C1").Value = "Food"
$ws.Range("D1").Value = "Calories"
$ws.Range("E1").Value = "Genotype"
$ws.Range("F1").Value = "BMI"

$ws.Range("A2").Value = "Leptin"
$ws.Range("A3").Value = "Food"
$ws.Range("A4").Value = "Calories"
$ws.Range("A5").Value = "Genotype"
$ws.Range("A6").Value = "BMI"

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0.068358738951585313
$ws.Range("D2").Value = 0.11774217370454802
$ws.Range("E2").Value = 0.13758362776769406
$ws.Range("F2").Value = 0.52721724605170217

$ws.Range("B3").Value = 0.068358738951585313
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.064616787549566707
$ws.Range("E3").Value = 0.33759358285959173
$ws.Range("F3").Value = 0.49226443591609892

$ws.Range("B4").Value = 0.11774217370454802
$ws.Range("C4").Value = 0.064616787549566707
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0.17529045003438418
$ws.Range("F4").Value = 0.51720675006084749

$ws.Range("B5").Value = 0.13758362776769406
$ws.Range("C5").Value = 0.33759358285959173
$ws.Range("D5").Value = 0.17529045003438418
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0.64125442193753468

$ws.Range("B6").Value = 0.52721724605170217
$ws.Range("C6").Value = 0.49226443591609892
$ws.Range("D6").Value = 0.51720675006084749
$ws.Range("E6").Value = 0.64125442193753468
$ws.Range("F6").Value = 0
